$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that originally sat right
#    before "Starting Files: Lab 5 Part 2".
# ------------------------------------------------------------------
$hasGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hasGoBack) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. "Add validation messages for the movie title." -> bold "movie title"
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Add validation messages for the movie title*") {
        $r = $p.Range.Duplicate
        $found = $r.Find.Execute("movie title", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $r.Bold = 1
    }
}

# ------------------------------------------------------------------
# 3. "Add validation messages for the director." -> bold "director"
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Add validation messages for the director*") {
        $r = $p.Range.Duplicate
        $found = $r.Find.Execute("director", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $r.Bold = 1
    }
}

# ------------------------------------------------------------------
# 4. Re-add the "_GoBack" bookmark right after "Run the application"
#    (collapsed bookmark at the very end of that paragraph's text).
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Trim() -eq "Run the application") {
        $r = $p.Range.Duplicate
        $moved = $r.MoveEnd(1, -1) # exclude the paragraph mark
        $r.Collapse(0)             # collapse to the very end of the text
        $r.InsertAfter("X")        # temporary marker so the collapse point
                                    # is unambiguously *after* the real text
        $r.Collapse(1)             # collapse back to just before the marker
        $bm = $d.Bookmarks.Add("_GoBack", $r)
        $delRange = $d.Range($r.End, $r.End + 1)
        $delRange.Delete()         # remove the temporary marker
    }
}

# ------------------------------------------------------------------
# 5. Split the `("movieHunter", ["` run so "movieHunter" is its own run
#    (mirrors the surrounding proofed-word runs like "ngRoute").
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*angular.module*") {
        $r = $p.Range.Duplicate
        $found = $r.Find.Execute("movieHunter", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $r.Font.Name = "Lucida Console"
        $r.Font.NameBi = "Lucida Console"
    }
}
